$d = $word.ActiveDocument

# Helper: find the index (1-based) of the first paragraph whose text
# contains the given substring.
function Get-ParaIndexByText($doc, $searchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$searchText*") {
            return $i
        }
    }
    return -1
}

# Locate the two footer paragraphs that must be removed entirely:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#    pages. Original theme under Creative Commons Attribution"
# These two paragraphs sit between two blank paragraphs (one right after
# "LOB1018: Física I (Requisito)" and one right before the page-break
# paragraph). Removing only the two text paragraphs leaves both blank
# paragraphs behind, so also collapse the now-redundant blank paragraph
# that used to directly precede "Ver no Jupiter...".
$idxVerNoJupiter = Get-ParaIndexByText $d "Ver no Jupiter Salvar em pdf Salvar em docx"
$idxCopyright = Get-ParaIndexByText $d "luizeleno@usp.br"
$idxBlank = $idxVerNoJupiter - 1

# Delete from the highest paragraph index to the lowest so earlier
# deletions never invalidate the indices computed above.
$d.Paragraphs.Item($idxCopyright).Range.Delete()
$d.Paragraphs.Item($idxVerNoJupiter).Range.Delete()
$d.Paragraphs.Item($idxBlank).Range.Delete()
